$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H300").Formula = "8098.5884"
$ws.Range("H301").Value2 = 8098.5884
$ws.Range("H302") = 8098.5884
